# Update cryptos worksheet values (prices, 1h volume %, and a few coin
# name/link/price row swaps) to reflect the refreshed GitHub Actions scrape.
#
# Numeric-looking price strings in column D are prefixed with a leading
# apostrophe so Excel stores them as literal text (preserving formatting
# such as trailing zeros and multi-dot separators) instead of silently
# re-parsing them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.556.78"
$ws.Range("E2").Value = "  +2.31%  "
$ws.Range("D3").Value = "2.723.57"
$ws.Range("E3").Value = "  +3.44%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'525.68"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("E6").Value = "  -1.24%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "'0.576"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").Value = "2.722.02"
$ws.Range("E9").Value = "  +2.40%  "
$ws.Range("E10").Value = "  +5.56%  "
$ws.Range("E11").Value = "  -0.18%  "
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("E13").Value = "  +3.05%  "
$ws.Range("D14").Value = "3.175.30"
$ws.Range("E14").Value = "  +2.43%  "
$ws.Range("D15").Value = "60.588.32"
$ws.Range("E15").Value = "  +2.43%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.850.22"
$ws.Range("E16").Value = "  +7.72%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'21.25"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "'345.67"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("E22").Value = "  +4.26%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'63.67"
$ws.Range("E24").Value = "  +2.98%  "
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("E26").Value = "  +3.80%  "
$ws.Range("D27").Value = "'0.994"
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  +0.85%  "
$ws.Range("D29").Value = "'7.28"
$ws.Range("E29").Value = "  +1.96%  "
$ws.Range("D30").Value = "'6.80"
$ws.Range("E30").Value = "  +8.20%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("D33").Value = "'19.05"
$ws.Range("E33").Value = "  +0.20%  "
$ws.Range("D34").Value = "'150.08"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "'4.25"
$ws.Range("E35").Value = "  +5.80%  "
$ws.Range("E36").Value = "  +8.24%  "
$ws.Range("D37").Value = "'0.941"
$ws.Range("E37").Value = "  -3.74%  "
$ws.Range("E38").Value = "  +6.64%  "
$ws.Range("D39").Value = "'0.871"
$ws.Range("E39").Value = "  +2.23%  "
$ws.Range("D40").Value = "'37.13"
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").Value = "'282.88"
$ws.Range("E42").Value = "  +1.04%  "
$ws.Range("D43").Value = "'20.10"
$ws.Range("E43").Value = "  +2.20%  "
$ws.Range("D44").Value = "'0.999"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "2.141.55"
$ws.Range("E45").Value = "  +7.41%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.0986"
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("E47").Value = "  -0.39%  "
$ws.Range("D48").Value = "'0.0538"
$ws.Range("E48").Value = "  +2.21%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'4.81"
$ws.Range("E49").Value = "  +0.95%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "'10.47"
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("E51").Value = "  +0.81%  "
